$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "runs", "balls", "fours" and "sixes" figures (columns C-F) for the two
# Tom Banton rows were swapped between row 2 and row 3. Keep the values
# stored as text (matching how the rest of the sheet stores its numbers)
# by applying a text number format before writing the new values.
$ws.Range("C2:F3").NumberFormat = "@"

$cols = @("C", "D", "E", "F")
foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $val2 = $cell2.Value2
    $val3 = $cell3.Value2

    $cell2.Value = $val3
    $cell3.Value = $val2
}
